$wb = $excel.ActiveWorkbook

$oldName = "Satisfaction between Age Groups"
$newName = "Performance between Age Groups"

$ws = $wb.Worksheets.Item($oldName)
$ws.Name = $newName

# Renaming a sheet via this automation surface does not cascade into chart
# series formulas the way interactive Excel does, so walk every chart on
# every worksheet and rewrite any SERIES() formula that still refers to the
# old sheet name.
foreach ($sheet in $wb.Worksheets) {
    foreach ($co in $sheet.ChartObjects()) {
        $chart = $co.Chart
        foreach ($series in $chart.SeriesCollection()) {
            if ($series.Formula -like "*$oldName*") {
                $series.Formula = $series.Formula.Replace($oldName, $newName)
            }
        }
    }
}
